$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.975.94"
$ws.Range("E2").Value = "  +4.50%  "

$ws.Range("D3").Value = "3.262.07"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'396.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'108.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("D7").Value = "'0.588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.12%  "

$ws.Range("D8").Value = "3.261.44"
$ws.Range("E8").Value = "  +2.99%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.627"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("D11").Value = "'39.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("D12").Value = "'0.0984"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.21%  "

$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").Value = "3.779.47"
$ws.Range("E14").Value = "  +2.78%  "

$ws.Range("D15").Value = "'8.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("D16").Value = "'19.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").Value = "3.262.83"
$ws.Range("E17").Value = "  +3.22%  "

$ws.Range("E18").Value = "  -2.47%  "

$ws.Range("E19").Value = "  +2.31%  "

$ws.Range("D20").Value = "56.888.57"
$ws.Range("E20").Value = "  +4.49%  "

$ws.Range("E21").Value = "  +1.76%  "

$ws.Range("E22").Value = "  +8.73%  "

$ws.Range("E23").Value = "  +1.31%  "

$ws.Range("D24").Value = "'296.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.78%  "

$ws.Range("E25").Value = "  +3.01%  "

$ws.Range("E26").Value = "  -2.48%  "

$ws.Range("E27").Value = "  +1.37%  "

$ws.Range("D28").Value = "'4.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.03%  "

$ws.Range("D29").Value = "'7.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("E30").Value = "  -4.06%  "

$ws.Range("D31").Value = "'0.168"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D33").Value = "'11.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.27%  "

$ws.Range("E34").Value = "  -2.83%  "

$ws.Range("D35").Value = "'40.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.57%  "

$ws.Range("E36").Value = "  -3.16%  "

$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("D38").Value = "'51.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").Value = "'3.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("D42").Value = "'138.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.85%  "

$ws.Range("E43").Value = "  +4.21%  "

$ws.Range("D44").Value = "'4.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "

$ws.Range("D45").Value = "'1.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("D46").Value = "'17.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("E47").Value = "  -3.24%  "

$ws.Range("D48").Value = "'22.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("E49").Value = "  +3.55%  "

$ws.Range("D50").Value = "2.168.77"
$ws.Range("E50").Value = "  +3.60%  "

$ws.Range("E51").Value = "  -5.71%  "
